# Bug 36035 70425 — Update translations of layouts and templates of Slide masters
# Applies German (de-DE) translations to the placeholder names/text that ship
# with the Slide Master, two of the Slide Layouts, and the Notes Master /
# Notes Page of the deck.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide Master (ppt/slideMasters/slideMaster1.xml)
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

# Shape 1 = Title placeholder ("Click to edit Master title style")
$masterTitle = $master.Shapes.Item(1)
$masterTitle.TextFrame.TextRange.Text = "Titelmasterformat durch Klicken bearbeiten"

# Shape 2 = Body placeholder (5 outline levels)
$masterBody = $master.Shapes.Item(2)
$masterBodyTr = $masterBody.TextFrame.TextRange
$masterBodyTr.Paragraphs(1).Text = "Textmasterformat bearbeiten"
$masterBodyTr.Paragraphs(2).Text = "Zweite Ebene"
$masterBodyTr.Paragraphs(3).Text = "Dritte Ebene"
$masterBodyTr.Paragraphs(4).Text = "Vierte Ebene"
$masterBodyTr.Paragraphs(5).Text = "Fünfte Ebene"

# ---------------------------------------------------------------------------
# 2) Slide Layout 5 ("Vergleich" / Comparison, ppt/slideLayouts/slideLayout5.xml)
#    Has two body placeholders, each with the same 5 outline levels.
# ---------------------------------------------------------------------------
$layouts = $master.CustomLayouts
$layoutCompare = $layouts.Item(5)

$compareBody1 = $layoutCompare.Shapes.Item(3)
$compareBody1Tr = $compareBody1.TextFrame.TextRange
$compareBody1Tr.Paragraphs(2).Text = "Zweite Ebene"
$compareBody1Tr.Paragraphs(3).Text = "Dritte Ebene"
$compareBody1Tr.Paragraphs(4).Text = "Vierte Ebene"
$compareBody1Tr.Paragraphs(5).Text = "Fünfte Ebene"

$compareBody2 = $layoutCompare.Shapes.Item(5)
$compareBody2Tr = $compareBody2.TextFrame.TextRange
$compareBody2Tr.Paragraphs(2).Text = "Zweite Ebene"
$compareBody2Tr.Paragraphs(3).Text = "Dritte Ebene"
$compareBody2Tr.Paragraphs(4).Text = "Vierte Ebene"
$compareBody2Tr.Paragraphs(5).Text = "Fünfte Ebene"

# ---------------------------------------------------------------------------
# 3) Slide Layout 9 ("Bild mit Überschrift" / Picture with Caption,
#    ppt/slideLayouts/slideLayout9.xml)
# ---------------------------------------------------------------------------
$layoutPicture = $layouts.Item(9)
$picturePlaceholder = $layoutPicture.Shapes.Item(2)
$picturePlaceholder.TextFrame.TextRange.Text = "Klicken Sie auf das Symbol, um ein Bild hinzuzufügen"

# ---------------------------------------------------------------------------
# 4) Notes Master (ppt/notesMasters/notesMaster1.xml)
#    Rename the placeholders and translate the body placeholder text.
# ---------------------------------------------------------------------------
$notesMaster = $p.NotesMaster

$notesMaster.Shapes.Item(1).Name = "Platzhalter für Überschrift 1"
$notesMaster.Shapes.Item(2).Name = "Datumsplatzhalter 2"
$notesMaster.Shapes.Item(3).Name = "Platzhalter für Folienbilder 3"
$notesMaster.Shapes.Item(4).Name = "Platzhalter für Notizen 4"
$notesMaster.Shapes.Item(5).Name = "Fußzeilenplatzhalter 5"
$notesMaster.Shapes.Item(6).Name = "Foliennummernplatzhalter 6"

$notesBody = $notesMaster.Shapes.Item(4)
$notesBody.TextFrame.TextRange.Text = "Textmasterformat bearbeiten" + [char]13 + "Zweite Ebene" + [char]13 + "Dritte Ebene" + [char]13 + "Vierte Ebene" + [char]13 + "Fünfte Ebene"

# ---------------------------------------------------------------------------
# 5) Notes Page of Slide 1 (ppt/notesSlides/notesSlide1.xml)
#    Rename the placeholders to their German equivalents.
# ---------------------------------------------------------------------------
$notesPage = $p.Slides.Item(1).NotesPage
$notesPage.Shapes.Item(1).Name = "Platzhalter für Folienbild 1"
$notesPage.Shapes.Item(2).Name = "Platzhalter für Notizen 2"
$notesPage.Shapes.Item(3).Name = "Foliennummernplatzhalter 3"
